# Add the two new "/add_vendor" entries to the "articels" sheet:
#   - row 4 (was an empty stub row) gets "/add_vendor"
#   - a brand-new row 9 gets "/add_vendor"
# and extend the sheet's used range down to row 99 (rows 10-99 stay blank,
# mirroring how the workbook looked after the upload: a big block of
# pre-formatted but empty rows below the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articels")

$ws.Range("A4").Value = "/add_vendor"
$ws.Range("A9").Value = "/add_vendor"

# Touch row 99 with a (no-op) formatting change so the sheet's dimension /
# used-range grows to include it, then group/ungroup rows 10-99 so every
# row in between materialises as a blank row in the saved sheet.
$ws.Range("A99").Borders.LineStyle = -4142

$ws.Rows("10:99").Group()
$ws.Rows("10:99").Ungroup()
